# Updates coin price/volume data to reflect the latest market snapshot.
# Also swaps the BinanceUSD/Polygon rows (9 and 10) which changed rank order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.852.92"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "1.708.69"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").Value = "'315.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "'0.4015"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.00%  "
$ws.Range("D8").Value = "'0.4049"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("B9").Value = "Polygon"
$ws.Range("C9").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D9").Value = "'1.476"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("B10").Value = "BinanceUSD"
$ws.Range("C10").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D10").Value = "'1.002"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("D11").Value = "'53.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").Value = "'0.08818"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").Value = "'26.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.70%  "
$ws.Range("D14").Value = "'7.536"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "'8.034"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "'0.00001345"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "1.710.49"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").Value = "'95.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.86%  "
$ws.Range("D19").Value = "'0.07171"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").Value = "'20.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.43%  "
$ws.Range("D21").Value = "'7.300"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").Value = "'14.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("D24").Value = "24.833.55"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").Value = "'2.359"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").Value = "'2.906"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("D27").Value = "'23.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").Value = "'6.136"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +17.61%  "
$ws.Range("D29").Value = "'161.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "'144.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.47%  "
$ws.Range("D31").Value = "'8.240"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.69%  "
$ws.Range("D32").Value = "'2.280"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +14.93%  "
$ws.Range("D33").Value = "1.912.06"
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("D34").Value = "'0.08677"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.07%  "
$ws.Range("D35").Value = "'0.03207"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.70%  "
$ws.Range("D36").Value = "'7.297"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("D37").Value = "'1.032"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("E38").Value = "  +4.46%  "
$ws.Range("E39").Value = "  +7.61%  "
$ws.Range("D40").Value = "'0.09482"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.72%  "
$ws.Range("D41").Value = "'10.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").Value = "'14.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("D43").Value = "'1.485"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("D44").Value = "'17.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.60%  "
$ws.Range("D45").Value = "'2.723"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.12%  "
$ws.Range("D46").Value = "'0.7449"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.61%  "
$ws.Range("D47").Value = "'4.218"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").Value = "'1.382"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.43%  "
$ws.Range("D49").Value = "'1.002"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "'140.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").Value = "'0.08403"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.50%  "
